$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price/Volume columns keep their original text (string) cell type,
# since several values look like plain numbers (e.g. "213.44", "0.995")
# and would otherwise be auto-converted to numeric cells by Excel.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range('D2').Value = '29.829.82'
$ws.Range('E2').Value = '  +0.83%  '
$ws.Range('D3').Value = '1.622.51'
$ws.Range('E3').Value = '  +1.01%  '
$ws.Range('E4').Value = '  -0.31%  '
$ws.Range('D5').Value = '213.44'
$ws.Range('E5').Value = '  +0.38%  '
$ws.Range('E6').Value = '  -0.73%  '
$ws.Range('D7').Value = '0.995'
$ws.Range('E7').Value = '  -0.31%  '
$ws.Range('D8').Value = '29.37'
$ws.Range('E8').Value = '  +9.00%  '
$ws.Range('E9').Value = '  +3.16%  '
$ws.Range('E10').Value = '  +0.95%  '
$ws.Range('E11').Value = '  +0.29%  '
$ws.Range('D12').Value = '1.855.45'
$ws.Range('E12').Value = '  +1.12%  '
$ws.Range('D13').Value = '1.623.32'
$ws.Range('E13').Value = '  +1.10%  '
$ws.Range('E14').Value = '  +5.49%  '
$ws.Range('D15').Value = '3.90'
$ws.Range('E15').Value = '  +4.88%  '
$ws.Range('D16').Value = '29.876.65'
$ws.Range('E16').Value = '  +1.00%  '
$ws.Range('E17').Value = '  +15.86%  '
$ws.Range('E18').Value = '  +1.34%  '
$ws.Range('D19').Value = '242.67'
$ws.Range('E19').Value = '  +0.83%  '
$ws.Range('E20').Value = '  +1.92%  '
$ws.Range('E21').Value = '  -0.26%  '
$ws.Range('D22').Value = '4.10'
$ws.Range('E22').Value = '  +2.96%  '
$ws.Range('D23').Value = '9.57'
$ws.Range('E23').Value = '  +3.73%  '
$ws.Range('E24').Value = '  +2.34%  '
$ws.Range('D25').Value = '156.41'
$ws.Range('E25').Value = '  +1.28%  '
$ws.Range('D26').Value = '15.60'
$ws.Range('E26').Value = '  +2.23%  '
$ws.Range('E27').Value = '  +1.34%  '
$ws.Range('E28').Value = '  +2.71%  '
$ws.Range('D29').Value = '0.996'
$ws.Range('E29').Value = '  -0.25%  '
$ws.Range('D30').Value = '0.0487'
$ws.Range('E30').Value = '  +3.23%  '
$ws.Range('E31').Value = '  +4.84%  '
$ws.Range('E32').Value = '  +3.03%  '
$ws.Range('E33').Value = '  +3.70%  '
$ws.Range('D34').Value = '1.424.64'
$ws.Range('E34').Value = '  +1.10%  '
$ws.Range('E35').Value = '  +6.51%  '
$ws.Range('E36').Value = '  -0.60%  '
$ws.Range('D37').Value = '2.86'
$ws.Range('E37').Value = '  +1.63%  '
$ws.Range('E38').Value = '  -0.62%  '
$ws.Range('E39').Value = '  +2.81%  '
$ws.Range('D40').Value = '0.554'
$ws.Range('E40').Value = '  +2.97%  '
$ws.Range('E41').Value = '  +3.11%  '
$ws.Range('D42').Value = '0.830'
$ws.Range('E42').Value = '  +3.84%  '
$ws.Range('E43').Value = '  -0.23%  '
$ws.Range('D44').Value = '54.27'
$ws.Range('E44').Value = '  +1.21%  '
$ws.Range('D45').Value = '68.99'
$ws.Range('E45').Value = '  +4.71%  '
$ws.Range('E46').Value = '  +19.26%  '
$ws.Range('D47').Value = '0.994'
$ws.Range('E47').Value = '  -0.37%  '
$ws.Range('D48').Value = '5.41'
$ws.Range('E48').Value = '  +2.39%  '
$ws.Range('D49').Value = '1.763.94'
$ws.Range('E49').Value = '  +1.02%  '
$ws.Range('D50').Value = '88.01'
$ws.Range('E50').Value = '  +1.51%  '
$ws.Range('E51').Value = '  +12.35%  '
